$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# --- Text / status / timestamp updates (Status moved from "In Translation" to "Ready for handoff") ---

# Overview sheet: row 2 holds the per-language status (E/F) and the "Latest HO Xliff
# Generate Date" (G), which moves forward a few seconds to reflect the new handoff.
$overview.Range("E2").Value = "Ready for handoff"
$overview.Range("F2").Value = "Ready for handoff"
$overview.Range("G2").Value = "2016-09-03 01:03:07"

# zh-cn sheet: Status column (C) + Latest Handoff Datetime (H)
$zhcn.Range("C2").Value = "Ready for handoff"
$zhcn.Range("H2").Value = "2016-09-03 01:02:59"

# de-de sheet: Status column (C); its "Latest Handoff Datetime" (H) mirrors the Overview date
$dede.Range("C2").Value = "Ready for handoff"
$dede.Range("H2").Value = "2016-09-03 01:03:07"

# --- Column width updates (status columns widened to fit "Ready for handoff") ---
# ColumnWidth is expressed in characters; the previous width corresponded to ~12.5
# characters, the new one to ~16.3 characters (wide enough for the longer text).
$overview.Columns.Item(5).ColumnWidth = 16.3
$overview.Columns.Item(6).ColumnWidth = 16.3
$zhcn.Columns.Item(3).ColumnWidth = 16.3
$dede.Columns.Item(3).ColumnWidth = 16.3
